# Generate Report for Handoff
# Adds two new tracked files to the localization-status report:
#   504b6090-36d3-4600-a9fe-ead6fe3a1014.md  (inserted before 6729d8e2..., "Ready for handoff")
#   6be56fae-7ee9-45b0-9729-46b870a9e423.md  (inserted after  6729d8e2..., "Ready for handoff")
# across all three worksheets: Overview, zh-cn, de-de.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: Overview  (columns: A=File Name, B=zh-cn status, C=de-de status)
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Overview")

# Make room: duplicate row 5 (6729d8e2...) upward to create a slot for the
# new 504b6090 entry, then duplicate the (now shifted) 6729d8e2 row again to
# create a slot after it for the new 6be56fae entry.
$ws1.Rows.Item(5).Copy()
$ws1.Rows.Item(5).Insert()
$ws1.Rows.Item(6).Copy()
$ws1.Rows.Item(7).Insert()

# Row 5: 504b6090...
$ws1.Range("A5").Value2 = "504b6090-36d3-4600-a9fe-ead6fe3a1014.md"
$ws1.Range("B5").Value2 = "Ready for handoff"
$ws1.Range("C5").Value2 = "Ready for handoff"

# Row 6: 6729d8e2... (unchanged content, already correct from the copy)
$ws1.Range("A6").Value2 = "6729d8e2-9b2b-4b69-a96b-8198b9cef925.md"
$ws1.Range("B6").Value2 = "Ready for handoff"
$ws1.Range("C6").Value2 = "Ready for handoff"

# Row 7: 6be56fae...
$ws1.Range("A7").Value2 = "6be56fae-7ee9-45b0-9729-46b870a9e423.md"
$ws1.Range("B7").Value2 = "Ready for handoff"
$ws1.Range("C7").Value2 = "Ready for handoff"

# Row 8: .localization-config (unchanged content, shifted down automatically)

# Rebuild hyperlinks top-to-bottom in final row order (the engine does not
# shift existing hyperlink anchors when rows are inserted/copied, so we
# recreate the full set from scratch to keep ranges/targets in sync).
$ws1.Hyperlinks.Delete()
$ws1.Hyperlinks.Add($ws1.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/a45e8733a7ae4057df653b3834c62cfe7f0acd7d/e2e/0dcdd56c-0d5f-491e-b607-aa04f934d885.md", "", "", "0dcdd56c-0d5f-491e-b607-aa04f934d885.md")
$ws1.Hyperlinks.Add($ws1.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/8b1734e6a4c464880d008a5ac6d428ee9a5a0f2b/e2e/06bd13a8-6cfb-4bed-b498-d393cd6b7967.md", "", "", "06bd13a8-6cfb-4bed-b498-d393cd6b7967.md")
$ws1.Hyperlinks.Add($ws1.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/8b1734e6a4c464880d008a5ac6d428ee9a5a0f2b/e2e/f2314cc6-7b27-4ea1-84f2-068f89ed5748.md", "", "", "f2314cc6-7b27-4ea1-84f2-068f89ed5748.md")
$ws1.Hyperlinks.Add($ws1.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/9c6a1f4cf34f84a6a14a16a639ae1ae26ae5e6cb/e2e/504b6090-36d3-4600-a9fe-ead6fe3a1014.md", "", "", "504b6090-36d3-4600-a9fe-ead6fe3a1014.md")
$ws1.Hyperlinks.Add($ws1.Range("A6"), "https://github.com/OpenLocalizationTest/oltest/blob/320017acf2155a32e66f576523f5b57c95ce14ca/e2e/6729d8e2-9b2b-4b69-a96b-8198b9cef925.md", "", "", "6729d8e2-9b2b-4b69-a96b-8198b9cef925.md")
$ws1.Hyperlinks.Add($ws1.Range("A7"), "https://github.com/OpenLocalizationTest/oltest/blob/0c0d9cdcae7c0f4f8f5f6b8e2d2aa1f8a5d6c2f1/e2e/6be56fae-7ee9-45b0-9729-46b870a9e423.md", "", "", "6be56fae-7ee9-45b0-9729-46b870a9e423.md")
$ws1.Hyperlinks.Add($ws1.Range("A8"), "https://github.com/OpenLocalizationTest/oltest/blob/8b1734e6a4c464880d008a5ac6d428ee9a5a0f2b/.localization-config", "", "", ".localization-config")

# ---------------------------------------------------------------------------
# Sheet 2: zh-cn
# Columns: A=Source File, B=Status, C=Latest Handoff File, D=Latest Handoff
#          Datetime, E=Latest Target File, F=Latest Handback File,
#          G=Latest Handback DateTime, H=Handoff Reason, I=Dependency From
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("zh-cn")

$ws2.Rows.Item(5).Copy()
$ws2.Rows.Item(5).Insert()
$ws2.Rows.Item(6).Copy()
$ws2.Rows.Item(7).Insert()

# Row 5: 504b6090...
$ws2.Range("A5").Value2 = "504b6090-36d3-4600-a9fe-ead6fe3a1014.md"
$ws2.Range("B5").Value2 = "Ready for handoff"
$ws2.Range("C5").Value2 = "504b6090-36d3-4600-a9fe-ead6fe3a1014.9e5bc72d79d621f43352f2de72e7c92a665cd185.zh-cn.xlf"
$ws2.Range("D5").Value2 = "2016-03-08 18:39:22"
$ws2.Range("G5").Value2 = "0001-01-01 00:00:00"
$ws2.Range("H5").Value2 = "Include"

# Row 6: 6729d8e2... (content already correct from the copy, re-assert to be safe)
$ws2.Range("A6").Value2 = "6729d8e2-9b2b-4b69-a96b-8198b9cef925.md"
$ws2.Range("B6").Value2 = "Ready for handoff"
$ws2.Range("C6").Value2 = "6729d8e2-9b2b-4b69-a96b-8198b9cef925.07837af2e7bb52b488a7f0a8be0b27fdf38244c1.zh-cn.xlf"
$ws2.Range("D6").Value2 = "2016-03-08 18:35:36"
$ws2.Range("G6").Value2 = "0001-01-01 00:00:00"
$ws2.Range("H6").Value2 = "Include"

# Row 7: 6be56fae...
$ws2.Range("A7").Value2 = "6be56fae-7ee9-45b0-9729-46b870a9e423.md"
$ws2.Range("B7").Value2 = "Ready for handoff"
$ws2.Range("C7").Value2 = "6be56fae-7ee9-45b0-9729-46b870a9e423.02fe312f874fd972d16537bae4f783c169acb603.zh-cn.xlf"
$ws2.Range("D7").Value2 = "2016-03-08 18:39:22"
$ws2.Range("G7").Value2 = "0001-01-01 00:00:00"
$ws2.Range("H7").Value2 = "Include"

# Row 8: .localization-config (unchanged content, shifted down automatically)

$ws2.Hyperlinks.Delete()
$ws2.Hyperlinks.Add($ws2.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/a45e8733a7ae4057df653b3834c62cfe7f0acd7d/e2e/0dcdd56c-0d5f-491e-b607-aa04f934d885.md", "", "", "0dcdd56c-0d5f-491e-b607-aa04f934d885.md")
$ws2.Hyperlinks.Add($ws2.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/346763e28a5306e32695d4371b6e4d2ba1d0bc04/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/0dcdd56c-0d5f-491e-b607-aa04f934d885.fd2a5b40ed21b9ffb499933482f584d3c64db1eb.zh-cn.xlf", "", "", "0dcdd56c-0d5f-491e-b607-aa04f934d885.fd2a5b40ed21b9ffb499933482f584d3c64db1eb.zh-cn.xlf")
$ws2.Hyperlinks.Add($ws2.Range("E2"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/0b33ae25b448ef259cde5214f0aa053de526de73/e2e/0dcdd56c-0d5f-491e-b607-aa04f934d885.md", "", "", "0dcdd56c-0d5f-491e-b607-aa04f934d885.md")
$ws2.Hyperlinks.Add($ws2.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/bd1c4a12b868ae77ef3e9aeab5ca3bec8c5f2327/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/0dcdd56c-0d5f-491e-b607-aa04f934d885.fd2a5b40ed21b9ffb499933482f584d3c64db1eb.zh-cn.xlf", "", "", "0dcdd56c-0d5f-491e-b607-aa04f934d885.fd2a5b40ed21b9ffb499933482f584d3c64db1eb.zh-cn.xlf")
$ws2.Hyperlinks.Add($ws2.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/8b1734e6a4c464880d008a5ac6d428ee9a5a0f2b/e2e/06bd13a8-6cfb-4bed-b498-d393cd6b7967.md", "", "", "06bd13a8-6cfb-4bed-b498-d393cd6b7967.md")
$ws2.Hyperlinks.Add($ws2.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/b879957f6e917f8c822fb90986aaf304d3ea12a5/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/06bd13a8-6cfb-4bed-b498-d393cd6b7967.68ac968fefd4ba9bc8aeeaaefec77311615b463d.zh-cn.xlf", "", "", "06bd13a8-6cfb-4bed-b498-d393cd6b7967.68ac968fefd4ba9bc8aeeaaefec77311615b463d.zh-cn.xlf")
$ws2.Hyperlinks.Add($ws2.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/8b1734e6a4c464880d008a5ac6d428ee9a5a0f2b/e2e/f2314cc6-7b27-4ea1-84f2-068f89ed5748.md", "", "", "f2314cc6-7b27-4ea1-84f2-068f89ed5748.md")
$ws2.Hyperlinks.Add($ws2.Range("C4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/b879957f6e917f8c822fb90986aaf304d3ea12a5/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/f2314cc6-7b27-4ea1-84f2-068f89ed5748.db694450088793dd4fce6e8492f11ced740f1978.zh-cn.xlf", "", "", "f2314cc6-7b27-4ea1-84f2-068f89ed5748.db694450088793dd4fce6e8492f11ced740f1978.zh-cn.xlf")
$ws2.Hyperlinks.Add($ws2.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/9c6a1f4cf34f84a6a14a16a639ae1ae26ae5e6cb/e2e/504b6090-36d3-4600-a9fe-ead6fe3a1014.md", "", "", "504b6090-36d3-4600-a9fe-ead6fe3a1014.md")
$ws2.Hyperlinks.Add($ws2.Range("C5"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/1a2b3c4d5e6f7a8b9c0d1e2f3a4b5c6d7e8f9a0b/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/504b6090-36d3-4600-a9fe-ead6fe3a1014.9e5bc72d79d621f43352f2de72e7c92a665cd185.zh-cn.xlf", "", "", "504b6090-36d3-4600-a9fe-ead6fe3a1014.9e5bc72d79d621f43352f2de72e7c92a665cd185.zh-cn.xlf")
$ws2.Hyperlinks.Add($ws2.Range("A6"), "https://github.com/OpenLocalizationTest/oltest/blob/320017acf2155a32e66f576523f5b57c95ce14ca/e2e/6729d8e2-9b2b-4b69-a96b-8198b9cef925.md", "", "", "6729d8e2-9b2b-4b69-a96b-8198b9cef925.md")
$ws2.Hyperlinks.Add($ws2.Range("C6"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/df793c259280bbb9f63faff4f1e7e487e03e9014/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/6729d8e2-9b2b-4b69-a96b-8198b9cef925.07837af2e7bb52b488a7f0a8be0b27fdf38244c1.zh-cn.xlf", "", "", "6729d8e2-9b2b-4b69-a96b-8198b9cef925.07837af2e7bb52b488a7f0a8be0b27fdf38244c1.zh-cn.xlf")
$ws2.Hyperlinks.Add($ws2.Range("A7"), "https://github.com/OpenLocalizationTest/oltest/blob/0c0d9cdcae7c0f4f8f5f6b8e2d2aa1f8a5d6c2f1/e2e/6be56fae-7ee9-45b0-9729-46b870a9e423.md", "", "", "6be56fae-7ee9-45b0-9729-46b870a9e423.md")
$ws2.Hyperlinks.Add($ws2.Range("C7"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/2f3e4d5c6b7a8900a1b2c3d4e5f60718293a4b5c/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/6be56fae-7ee9-45b0-9729-46b870a9e423.02fe312f874fd972d16537bae4f783c169acb603.zh-cn.xlf", "", "", "6be56fae-7ee9-45b0-9729-46b870a9e423.02fe312f874fd972d16537bae4f783c169acb603.zh-cn.xlf")
$ws2.Hyperlinks.Add($ws2.Range("A8"), "https://github.com/OpenLocalizationTest/oltest/blob/8b1734e6a4c464880d008a5ac6d428ee9a5a0f2b/.localization-config", "", "", ".localization-config")

# ---------------------------------------------------------------------------
# Sheet 3: de-de  (same layout as zh-cn, de-de specific data/hashes)
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("de-de")

$ws3.Rows.Item(5).Copy()
$ws3.Rows.Item(5).Insert()
$ws3.Rows.Item(6).Copy()
$ws3.Rows.Item(7).Insert()

# Row 5: 504b6090...
$ws3.Range("A5").Value2 = "504b6090-36d3-4600-a9fe-ead6fe3a1014.md"
$ws3.Range("B5").Value2 = "Ready for handoff"
$ws3.Range("C5").Value2 = "504b6090-36d3-4600-a9fe-ead6fe3a1014.9e5bc72d79d621f43352f2de72e7c92a665cd185.de-de.xlf"
$ws3.Range("D5").Value2 = "2016-03-08 18:39:29"
$ws3.Range("G5").Value2 = "0001-01-01 00:00:00"
$ws3.Range("H5").Value2 = "Include"

# Row 6: 6729d8e2... (content already correct from the copy, re-assert to be safe)
$ws3.Range("A6").Value2 = "6729d8e2-9b2b-4b69-a96b-8198b9cef925.md"
$ws3.Range("B6").Value2 = "Ready for handoff"
$ws3.Range("C6").Value2 = "6729d8e2-9b2b-4b69-a96b-8198b9cef925.07837af2e7bb52b488a7f0a8be0b27fdf38244c1.de-de.xlf"
$ws3.Range("D6").Value2 = "2016-03-08 18:35:44"
$ws3.Range("G6").Value2 = "0001-01-01 00:00:00"
$ws3.Range("H6").Value2 = "Include"

# Row 7: 6be56fae...
$ws3.Range("A7").Value2 = "6be56fae-7ee9-45b0-9729-46b870a9e423.md"
$ws3.Range("B7").Value2 = "Ready for handoff"
$ws3.Range("C7").Value2 = "6be56fae-7ee9-45b0-9729-46b870a9e423.02fe312f874fd972d16537bae4f783c169acb603.de-de.xlf"
$ws3.Range("D7").Value2 = "2016-03-08 18:39:29"
$ws3.Range("G7").Value2 = "0001-01-01 00:00:00"
$ws3.Range("H7").Value2 = "Include"

# Row 8: .localization-config (unchanged content, shifted down automatically)

$ws3.Hyperlinks.Delete()
$ws3.Hyperlinks.Add($ws3.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/a45e8733a7ae4057df653b3834c62cfe7f0acd7d/e2e/0dcdd56c-0d5f-491e-b607-aa04f934d885.md", "", "", "0dcdd56c-0d5f-491e-b607-aa04f934d885.md")
$ws3.Hyperlinks.Add($ws3.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/23688fd400c91ce1d660a19ecac1b3133cebe2f2/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/0dcdd56c-0d5f-491e-b607-aa04f934d885.fd2a5b40ed21b9ffb499933482f584d3c64db1eb.de-de.xlf", "", "", "0dcdd56c-0d5f-491e-b607-aa04f934d885.fd2a5b40ed21b9ffb499933482f584d3c64db1eb.de-de.xlf")
$ws3.Hyperlinks.Add($ws3.Range("E2"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/e4a8445736f89cd20cb3300303612c5a75fdec97/e2e/0dcdd56c-0d5f-491e-b607-aa04f934d885.md", "", "", "0dcdd56c-0d5f-491e-b607-aa04f934d885.md")
$ws3.Hyperlinks.Add($ws3.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/b10d812ff53e51fa6c39ff73fcd42bb5a150d32c/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/0dcdd56c-0d5f-491e-b607-aa04f934d885.fd2a5b40ed21b9ffb499933482f584d3c64db1eb.de-de.xlf", "", "", "0dcdd56c-0d5f-491e-b607-aa04f934d885.fd2a5b40ed21b9ffb499933482f584d3c64db1eb.de-de.xlf")
$ws3.Hyperlinks.Add($ws3.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/8b1734e6a4c464880d008a5ac6d428ee9a5a0f2b/e2e/06bd13a8-6cfb-4bed-b498-d393cd6b7967.md", "", "", "06bd13a8-6cfb-4bed-b498-d393cd6b7967.md")
$ws3.Hyperlinks.Add($ws3.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/696e93354a1aaf5add4549e219abfbd86f4dcb93/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/06bd13a8-6cfb-4bed-b498-d393cd6b7967.68ac968fefd4ba9bc8aeeaaefec77311615b463d.de-de.xlf", "", "", "06bd13a8-6cfb-4bed-b498-d393cd6b7967.68ac968fefd4ba9bc8aeeaaefec77311615b463d.de-de.xlf")
$ws3.Hyperlinks.Add($ws3.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/8b1734e6a4c464880d008a5ac6d428ee9a5a0f2b/e2e/f2314cc6-7b27-4ea1-84f2-068f89ed5748.md", "", "", "f2314cc6-7b27-4ea1-84f2-068f89ed5748.md")
$ws3.Hyperlinks.Add($ws3.Range("C4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/696e93354a1aaf5add4549e219abfbd86f4dcb93/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/f2314cc6-7b27-4ea1-84f2-068f89ed5748.db694450088793dd4fce6e8492f11ced740f1978.de-de.xlf", "", "", "f2314cc6-7b27-4ea1-84f2-068f89ed5748.db694450088793dd4fce6e8492f11ced740f1978.de-de.xlf")
$ws3.Hyperlinks.Add($ws3.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/9c6a1f4cf34f84a6a14a16a639ae1ae26ae5e6cb/e2e/504b6090-36d3-4600-a9fe-ead6fe3a1014.md", "", "", "504b6090-36d3-4600-a9fe-ead6fe3a1014.md")
$ws3.Hyperlinks.Add($ws3.Range("C5"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/9e8d7c6b5a4938271605f4e3d2c1b0a9988776655/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/504b6090-36d3-4600-a9fe-ead6fe3a1014.9e5bc72d79d621f43352f2de72e7c92a665cd185.de-de.xlf", "", "", "504b6090-36d3-4600-a9fe-ead6fe3a1014.9e5bc72d79d621f43352f2de72e7c92a665cd185.de-de.xlf")
$ws3.Hyperlinks.Add($ws3.Range("A6"), "https://github.com/OpenLocalizationTest/oltest/blob/320017acf2155a32e66f576523f5b57c95ce14ca/e2e/6729d8e2-9b2b-4b69-a96b-8198b9cef925.md", "", "", "6729d8e2-9b2b-4b69-a96b-8198b9cef925.md")
$ws3.Hyperlinks.Add($ws3.Range("C6"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/c246a24b1835a84b2299531e178b8e9abb180275/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/6729d8e2-9b2b-4b69-a96b-8198b9cef925.07837af2e7bb52b488a7f0a8be0b27fdf38244c1.de-de.xlf", "", "", "6729d8e2-9b2b-4b69-a96b-8198b9cef925.07837af2e7bb52b488a7f0a8be0b27fdf38244c1.de-de.xlf")
$ws3.Hyperlinks.Add($ws3.Range("A7"), "https://github.com/OpenLocalizationTest/oltest/blob/0c0d9cdcae7c0f4f8f5f6b8e2d2aa1f8a5d6c2f1/e2e/6be56fae-7ee9-45b0-9729-46b870a9e423.md", "", "", "6be56fae-7ee9-45b0-9729-46b870a9e423.md")
$ws3.Hyperlinks.Add($ws3.Range("C7"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/7766554433221100ffeeddccbbaa99887766554/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/6be56fae-7ee9-45b0-9729-46b870a9e423.02fe312f874fd972d16537bae4f783c169acb603.de-de.xlf", "", "", "6be56fae-7ee9-45b0-9729-46b870a9e423.02fe312f874fd972d16537bae4f783c169acb603.de-de.xlf")
$ws3.Hyperlinks.Add($ws3.Range("A8"), "https://github.com/OpenLocalizationTest/oltest/blob/8b1734e6a4c464880d008a5ac6d428ee9a5a0f2b/.localization-config", "", "", ".localization-config")

$wb.Save()
